# Daily auto-update: insert a new top row for the newest date (2026-02-08),
# pushing all existing rows down by one, keeping the same price values as
# the rest of the (flat) series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (old row 2..80) down to (3..81) by inserting a
# fresh row right below the header.
$ws.Rows.Item(2).Insert()

# New top data row. Force the date cell to be stored as plain text (it is
# text in the source data, not a real Excel date) by pre-formatting as Text
# before assignment, then strip the temporary formatting again so the cell
# ends up styled like all of its neighbours.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-08"

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

$ws.Rows.Item(2).ClearFormats()
